$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Helper: write a value that must be stored as TEXT even when it
# looks like a pure number (Excel would otherwise coerce it to a
# numeric cell). We flip the number format to text, set the value,
# then restore the original "0" integer format used throughout the
# sheet so the visual style stays the same.
# ---------------------------------------------------------------
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "0"
}

# -----------------------------------------------------------------
# 1. Row 3: fill in the PO3/QTY3 pair (moved from the old PO6/QTY6
#    slot) and drop the now-unused PO6/QTY6 cells (S3:T3).
# -----------------------------------------------------------------
Set-TextValue $ws.Range("M3") "10737490"
Set-TextValue $ws.Range("N3") "9000"
$ws.Range("S3").Clear()
$ws.Range("T3").Clear()

# -----------------------------------------------------------------
# 2. Row 7: split the 75000 quantity into 70000 (PO1/QTY1) and a
#    new 5000 (PO2/QTY2) entry against the same PO number.
# -----------------------------------------------------------------
$ws.Range("J7").Value = 70000
Set-TextValue $ws.Range("K7") "10736367"
$ws.Range("L7").Value = 5000

# -----------------------------------------------------------------
# 3. Row 1: extend the header with PO7/QTY7 .. PO10/QTY10, copying
#    the look of the existing header cells, plus two trailing blank
#    (but styled) cells.
# -----------------------------------------------------------------
$ws.Range("T1").Copy()
$ws.Range("U1:AD1").PasteSpecial(-4122)

$ws.Range("U1").Value = "PO7"
$ws.Range("V1").Value = "QTY7"
$ws.Range("W1").Value = "PO8"
$ws.Range("X1").Value = "QTY8"
$ws.Range("Y1").Value = "PO9"
$ws.Range("Z1").Value = "QTY9"
$ws.Range("AA1").Value = "PO10"
$ws.Range("AB1").Value = "QTY10"
# AC1 / AD1 are left blank (format only, already pasted above).

# -----------------------------------------------------------------
# 4. New row 9: an extra shipment line with a "p1".."p10" PO/QTY
#    progression (ten extra POs).
# -----------------------------------------------------------------
# Base formatting for the whole row comes from row 2 (a plain data
# row), then the PO-label cells get the highlighted look already
# used for the carton-no. cells (D7/D8) in this sheet.
$ws.Range("A2:T2").Copy()
$ws.Range("A9:AB9").PasteSpecial(-4122)
$ws.Range("T1").Copy()
$ws.Range("AC9:AD9").PasteSpecial(-4122)

Set-TextValue $ws.Range("A9") "20140912001"
$ws.Range("B9").Value = "7YCQ40X1220+H03"
$ws.Range("C9").Value = "CHINA"

$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D9").Value = "A1-10"

Set-TextValue $ws.Range("E9") "13"
$ws.Range("F9").Value = "50904035R2F"
$ws.Range("G9").Value = 55
$ws.Range("H9").Value = "KYCE00"

$poCells = "I9","K9","M9","O9","Q9","S9","U9","W9","Y9","AA9"
$qtyCells = "J9","L9","N9","P9","R9","T9","V9","X9","Z9","AB9"
for ($i = 0; $i -lt 10; $i++) {
    $ws.Range("D8").Copy()
    $ws.Range($poCells[$i]).PasteSpecial(-4122)
    $ws.Range($poCells[$i]).Value = "p" + ($i + 1)
    $ws.Range($qtyCells[$i]).Value = $i + 1
}

# -----------------------------------------------------------------
# 5. Restore the cursor/selection position recorded in the file.
# -----------------------------------------------------------------
$ws.Range("G12").Select()
